$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (preserves original string formatting,
# e.g. trailing zeros / thousand-dot separators).
$textCells = @("D5", "D6", "D9", "D10", "D11", "D13", "D15", "D16", "D18", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D34", "D36", "D38", "D39", "D44", "D47", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '62.736.80'
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").Value = '3.466.74'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '414.61'
$ws.Range("E5").Value = '  +1.34%  '
$ws.Range("D6").Value = '130.82'
$ws.Range("E6").Value = '  +1.75%  '
$ws.Range("E7").Value = '  -1.40%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '0.726'
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("D10").Value = '0.142'
$ws.Range("E10").Value = '  +0.27%  '
$ws.Range("D11").Value = '42.70'
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("E12").Value = '  +6.44%  '
$ws.Range("D13").Value = '0.0000218'
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("D14").Value = '4.019.45'
$ws.Range("E14").Value = '  +1.43%  '
$ws.Range("D15").Value = '0.140'
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").Value = '20.53'
$ws.Range("E16").Value = '  -4.06%  '
$ws.Range("D17").Value = '3.491.07'
$ws.Range("E17").Value = '  +1.24%  '
$ws.Range("D18").Value = '12.75'
$ws.Range("E18").Value = '  +1.75%  '
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("D20").Value = '62.683.48'
$ws.Range("E20").Value = '  +1.12%  '
$ws.Range("D21").Value = '466.58'
$ws.Range("E21").Value = '  +3.24%  '
$ws.Range("D22").Value = '90.84'
$ws.Range("E22").Value = '  -1.66%  '
$ws.Range("E23").Value = '  +1.93%  '
$ws.Range("D24").Value = '13.37'
$ws.Range("E24").Value = '  +3.05%  '
$ws.Range("D25").Value = '10.64'
$ws.Range("E25").Value = '  +21.34%  '
$ws.Range("D26").Value = '3.31'
$ws.Range("E26").Value = '  +2.06%  '
$ws.Range("D27").Value = '33.33'
$ws.Range("E27").Value = '  +0.75%  '
$ws.Range("D28").Value = '4.80'
$ws.Range("E28").Value = '  +0.56%  '
$ws.Range("D29").Value = '7.60'
$ws.Range("E29").Value = '  -1.28%  '
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("E31").Value = '  -4.07%  '
$ws.Range("E32").Value = '  -2.19%  '
$ws.Range("E33").Value = '  -1.66%  '
$ws.Range("D34").Value = '41.00'
$ws.Range("E34").Value = '  -5.02%  '
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").Value = '58.58'
$ws.Range("E36").Value = '  +7.81%  '
$ws.Range("E37").Value = '  -2.40%  '
$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '3.08'
$ws.Range("E39").Value = '  +5.15%  '
$ws.Range("E40").Value = '  -0.80%  '
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("E42").Value = '  -0.83%  '
$ws.Range("E43").Value = '  +7.00%  '
$ws.Range("D44").Value = '145.98'
$ws.Range("E44").Value = '  +3.01%  '
$ws.Range("E45").Value = '  +2.29%  '
$ws.Range("E46").Value = '  +4.42%  '
$ws.Range("D47").Value = '2.42'
$ws.Range("E47").Value = '  +12.59%  '
$ws.Range("D48").Value = '0.0₃0564'
$ws.Range("E48").Value = '  +32.09%  '
$ws.Range("D49").Value = '16.39'
$ws.Range("E49").Value = '  -1.60%  '
$ws.Range("D50").Value = '22.37'
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("E51").Value = '  +1.19%  '
